$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 69
$ws.Range("H69").Value = 8724.583000000001
$ws.Range("J69").Value = 8724.583000000001
$ws.Range("L69").Value = 26173.749
$ws.Range("N69").Value = -27921.749
# Row 72
$ws.Range("H72").Value = 8724.583000000001
$ws.Range("J72").Value = 8724.583000000001
$ws.Range("L72").Value = 78521.247
$ws.Range("N72").Value = -87257.247
# Row 88
$ws.Range("H88").Value = 13914.417
$ws.Range("J88").Value = 17220.555
$ws.Range("L88").Value = 17220.555
$ws.Range("N88").Value = -18032.555
# Row 91
$ws.Range("H91").Value = 13914.417
$ws.Range("J91").Value = 17220.555
$ws.Range("L91").Value = 17220.555
$ws.Range("N91").Value = -20028.555
# Row 92
$ws.Range("H92").Value = 1604
$ws.Range("I92").Value = 344
$ws.Range("J92").Value = 2864
$ws.Range("K92").Value = 344
$ws.Range("L92").Value = 2864
$ws.Range("M92").Value = 904
$ws.Range("N92").Value = -5360
# Row 101
$ws.Range("H101").Value = 1462.4546
$ws.Range("I101").Value = 730.2857
$ws.Range("J101").Value = 2743.75
$ws.Range("K101").Value = 2190.8571
$ws.Range("L101").Value = 8231.25
$ws.Range("M101").Value = -568.8571000000002
$ws.Range("N101").Value = -11475.25
# Row 103
$ws.Range("H103").Value = 2188.5
$ws.Range("J103").Value = 2372.125
$ws.Range("L103").Value = 7116.375
$ws.Range("N103").Value = -8288.375
# Row 112
$ws.Range("H112").Value = 1795.4375
$ws.Range("I112").Value = 607.8
$ws.Range("J112").Value = 2335.2727
$ws.Range("K112").Value = 1823.4
$ws.Range("L112").Value = 7005.8181
$ws.Range("M112").Value = -715.3999999999999
$ws.Range("N112").Value = -9221.8181
# Row 114
$ws.Range("H114").Value = 59984
$ws.Range("J114").Value = 59984
$ws.Range("L114").Value = 59984
$ws.Range("N114").Value = -68662
# Row 116
$ws.Range("H116").Value = 11434.333
$ws.Range("I116").Value = 10400
$ws.Range("J116").Value = 13503
$ws.Range("K116").Value = 10400
$ws.Range("L116").Value = 13503
$ws.Range("M116").Value = -6958
$ws.Range("N116").Value = -20387
# Row 118
$ws.Range("H118").Value = 349.2857
$ws.Range("I118").Value = 190.83333
$ws.Range("J118").Value = 1300
$ws.Range("K118").Value = 572.49999
$ws.Range("L118").Value = 3900
$ws.Range("M118").Value = 1084.50001
$ws.Range("N118").Value = -7214
# Row 121
$ws.Range("H121").Value = 2318.8572
$ws.Range("J121").Value = 2318.8572
$ws.Range("L121").Value = 6956.571599999999
$ws.Range("N121").Value = -10450.5716
# Row 132
$ws.Range("H132").Value = 2399.8
$ws.Range("J132").Value = 5006
$ws.Range("L132").Value = 15018
$ws.Range("N132").Value = -20078
# Row 135
$ws.Range("H135").Value = 1115.7
$ws.Range("I135").Value = 878.9375
$ws.Range("J135").Value = 2062.75
$ws.Range("K135").Value = 7910.4375
$ws.Range("L135").Value = 18564.75
$ws.Range("M135").Value = -5375.4375
$ws.Range("N135").Value = -23634.75

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 24
$ws.Range("H24").Value = 27833.334
$ws.Range("J24").Value = 27833.334
$ws.Range("L24").Value = 27833.334
$ws.Range("N24").Value = -28581.334
# Row 32
$ws.Range("H32").Value = 3047.4517
$ws.Range("I32").Value = 2621.7705
$ws.Range("K32").Value = 2621.7705
$ws.Range("M32").Value = -2334.7705
# Row 74
$ws.Range("H74").Value = 15152775
$ws.Range("I74").Value = 19608884
$ws.Range("K74").Value = 19608884
$ws.Range("M74").Value = -19608010
# Row 77
$ws.Range("H77").Value = 15152775
$ws.Range("I77").Value = 19608884
$ws.Range("K77").Value = 98044420
$ws.Range("M77").Value = -98040052
# Row 93
$ws.Range("H93").Value = 28600
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 28600
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 28600
$ws.Range("M93").ClearContents()
$ws.Range("N93").Value = -33592
# Row 100
$ws.Range("H100").Value = 27833.334
$ws.Range("J100").Value = 27833.334
$ws.Range("L100").Value = 27833.334
$ws.Range("N100").Value = -29997.334
# Row 102
$ws.Range("H102").Value = 2016
$ws.Range("I102").Value = 2110.8
$ws.Range("J102").Value = 1700
$ws.Range("K102").Value = 2110.8
$ws.Range("L102").Value = 1700
$ws.Range("M102").Value = -488.8000000000002
$ws.Range("N102").Value = -4944

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 3988.889
$ws.Range("I86").Value = 1798.6666
$ws.Range("J86").Value = 5084
$ws.Range("K86").Value = 1798.6666
$ws.Range("L86").Value = 5084
$ws.Range("M86").Value = -675.6666
$ws.Range("N86").Value = -7330
# Row 89
$ws.Range("H89").Value = 3988.889
$ws.Range("I89").Value = 1798.6666
$ws.Range("J89").Value = 5084
$ws.Range("K89").Value = 8993.333000000001
$ws.Range("L89").Value = 25420
$ws.Range("M89").Value = -3377.333000000001
$ws.Range("N89").Value = -36652
# Row 93
$ws.Range("H93").Value = 40000
$ws.Range("J93").Value = 40000
$ws.Range("L93").Value = 40000
$ws.Range("N93").Value = -43744
# Row 94
$ws.Range("H94").Value = 962.7917
$ws.Range("I94").Value = 896.6
$ws.Range("K94").Value = 896.6
$ws.Range("M94").Value = -445.6
# Row 134
$ws.Range("H134").Value = 1301.5
$ws.Range("I134").Value = 1305.0526
$ws.Range("J134").Value = 1279
$ws.Range("K134").Value = 3915.1578
$ws.Range("L134").Value = 3837
$ws.Range("M134").Value = -1380.1578
$ws.Range("N134").Value = -8907

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 103
$ws.Range("H103").Value = 29219.25
$ws.Range("I103").Value = 20000
$ws.Range("J103").Value = 32292.334
$ws.Range("K103").Value = 20000
$ws.Range("L103").Value = 32292.334
$ws.Range("M103").Value = -18828
$ws.Range("N103").Value = -34636.334
# Row 107
$ws.Range("H107").Value = 945.3200000000001
$ws.Range("I107").Value = 837.17645
$ws.Range("K107").Value = 837.17645
$ws.Range("M107").Value = 1082.82355
# Row 114
$ws.Range("H114").Value = 24500
$ws.Range("J114").Value = 24500
$ws.Range("L114").Value = 24500
$ws.Range("N114").Value = -33178
# Row 132
$ws.Range("H132").Value = 10262.429
$ws.Range("I132").Value = 8803.833000000001
$ws.Range("K132").Value = 26411.499
$ws.Range("M132").Value = -23881.499

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 20
$ws.Range("H20").Value = 2724.5
$ws.Range("I20").Value = 4949
$ws.Range("J20").Value = 500
$ws.Range("K20").Value = 14847
$ws.Range("L20").Value = 1500
$ws.Range("M20").Value = -14620
$ws.Range("N20").Value = -1954
# Row 57
$ws.Range("H57").Value = 849.75
$ws.Range("I57").Value = 699.5
$ws.Range("K57").Value = 2098.5
$ws.Range("M57").Value = -1539.5
# Row 92
$ws.Range("H92").Value = 2100
$ws.Range("I92").Value = 650
$ws.Range("J92").Value = 5000
$ws.Range("K92").Value = 1950
$ws.Range("L92").Value = 15000
$ws.Range("M92").Value = -702
$ws.Range("N92").Value = -17496
# Row 131
$ws.Range("H131").Value = 26622346
$ws.Range("I131").Value = 20836168
$ws.Range("J131").Value = 32408522
$ws.Range("K131").Value = 62508504
$ws.Range("L131").Value = 97225566
$ws.Range("M131").Value = -62503464
$ws.Range("N131").Value = -97235646

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 5825.6313
$ws.Range("I80").Value = 5066
$ws.Range("J80").Value = 6176.231
$ws.Range("K80").Value = 5066
$ws.Range("L80").Value = 6176.231
$ws.Range("M80").Value = -4068
$ws.Range("N80").Value = -8172.231
# Row 83
$ws.Range("H83").Value = 5825.6313
$ws.Range("I83").Value = 5066
$ws.Range("J83").Value = 6176.231
$ws.Range("K83").Value = 25330
$ws.Range("L83").Value = 30881.155
$ws.Range("M83").Value = -20338
$ws.Range("N83").Value = -40865.155
# Row 92
$ws.Range("H92").Value = 21500
$ws.Range("J92").Value = 21500
$ws.Range("L92").Value = 21500
$ws.Range("N92").Value = -25244
# Row 101
$ws.Range("H101").Value = 26116.75
$ws.Range("J101").Value = 26116.75
$ws.Range("L101").Value = 26116.75
$ws.Range("N101").Value = -32606.75
# Row 107
$ws.Range("H107").Value = 1408.7059
$ws.Range("I107").Value = 682.8889
$ws.Range("J107").Value = 2225.25
$ws.Range("K107").Value = 682.8889
$ws.Range("L107").Value = 2225.25
$ws.Range("M107").Value = 1237.1111
$ws.Range("N107").Value = -6065.25
# Row 113
$ws.Range("H113").Value = 4250.7
$ws.Range("I113").Value = 3044.4
$ws.Range("K113").Value = 3044.4
$ws.Range("M113").Value = -874.4000000000001
# Row 132
$ws.Range("H132").Value = 68821.875
$ws.Range("I132").Value = 103372.8
$ws.Range("J132").Value = 11237
$ws.Range("K132").Value = 310118.4
$ws.Range("L132").Value = 33711
$ws.Range("M132").Value = -307588.4
$ws.Range("N132").Value = -38771

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 6063.0835
$ws.Range("I61").Value = 1292.3334
$ws.Range("J61").Value = 10833.833
$ws.Range("K61").Value = 1292.3334
$ws.Range("L61").Value = 10833.833
$ws.Range("M61").Value = -1090.3334
$ws.Range("N61").Value = -11237.833
# Row 93
$ws.Range("H93").Value = 4056.5
$ws.Range("I93").Value = 4010.8333
$ws.Range("J93").Value = 4102.1665
$ws.Range("K93").Value = 4010.8333
$ws.Range("L93").Value = 4102.1665
$ws.Range("M93").Value = -2762.8333
$ws.Range("N93").Value = -6598.1665
# Row 100
$ws.Range("H100").Value = 10487.667
$ws.Range("I100").Value = 9993
$ws.Range("J100").Value = 11477
$ws.Range("K100").Value = 9993
$ws.Range("L100").Value = 11477
$ws.Range("M100").Value = -9452
$ws.Range("N100").Value = -12559
# Row 113
$ws.Range("H113").Value = 6063.0835
$ws.Range("I113").Value = 1292.3334
$ws.Range("J113").Value = 10833.833
$ws.Range("K113").Value = 1292.3334
$ws.Range("L113").Value = 10833.833
$ws.Range("M113").Value = 877.6666
$ws.Range("N113").Value = -15173.833
# Row 132
$ws.Range("H132").Value = 4206.7715
$ws.Range("I132").Value = 2421.077
$ws.Range("K132").Value = 7263.231000000001
$ws.Range("M132").Value = -4733.231000000001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 42
$ws.Range("H42").Value = 31049
$ws.Range("J42").Value = 31049
$ws.Range("L42").Value = 31049
$ws.Range("N42").Value = -31805
# Row 107
$ws.Range("H107").Value = 426.16666
$ws.Range("I107").Value = 374
$ws.Range("K107").Value = 1122
$ws.Range("M107").Value = 798
# Row 120
$ws.Range("H120").Value = 74678.5
$ws.Range("J120").Value = 74678.5
$ws.Range("L120").Value = 74678.5
$ws.Range("N120").Value = -84354.5
# Row 132
$ws.Range("H132").Value = 6518.1577
$ws.Range("I132").Value = 5824.4443
$ws.Range("K132").Value = 17473.3329
$ws.Range("M132").Value = -14943.3329
